# "Role Refreshing Bugs fixed"
# Sprint row for 2013-W51 (row 13/14 of the burndown table) was missing the
# "Remaining work" entry and had a stale headcount. Fill those two inputs in
# back in - the Burndown/EVM formulas (F/G/H/I/J/K columns) recompute from
# them automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

# Remaining work for the row that was left blank.
$ws.Range("E13").Value = 12

# Role count correction for the following sprint row.
$ws.Range("B14").Value = 12

# Leave the sheet's selection where the author last left it.
$ws.Activate()
$ws.Range("D17").Select()
